# Update countries & provincias Spain
# Applies the 25-Abril-2020 10:22 data refresh to the "Pais" sheet:
#   - updates the "last updated" timestamp in A1
#   - refreshes per-country case/death statistics
#   - re-sorts two country pairs whose "Casos totales" changed rank order
#     (Filipinas overtakes Chequia; Afganistan overtakes Nueva Zelanda)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp header -------------------------------------------------
$ws.Cells.Item(1,1).Value = "Datos actualizados a 25 de Abril de 2020 a las 10:22"

# --- Row data (row, Pais, CasosTotales, NuevosCasos, CasosActivos, Recuperados, CasosCriticos, MuertesHoy, Muertes)
$rows = @(
    @(8,  "Alemania",             155054, 55,  109800, 39487, 2908, 7,  5767),
    @(33, "Polonia",               11067, 175,   2126,  8442,  160, 5,   499),
    @(35, "Rumania",               10417,   0,   2817,  7025,  241, 8,   575),
    @(44, "Filipinas",              7294, 102,    792,  6008,    1, 17,  494),
    @(45, "Chequia",                7273,   0,   2389,  4669,   71, 1,   215),
    @(46, "Australia",              6695,  20,   5372,  1243,   43, 1,    80),
    @(48, "Malasia",                5742,  51,   3762,  1882,   36, 2,    98),
    @(64, "Kazajistan",             2482,  66,    604,  1853,   31, 0,    25),
    @(74, "Afganistan",             1463, 112,    188,  1228,    7, 4,    47),
    @(75, "Nueva Zelanda",          1461,   5,   1118,   325,    1, 1,    18),
    @(76, "Camerun",                1430,   0,    668,   719,   20, 1,    41),
    @(77, "Lituania",               1426,  16,    460,   925,   17, 1,    41),
    @(78, "Bosnia y Herzegovina",   1421,   0,    538,   828,    4, 0,    55),
    @(79, "Eslovenia",              1373,   0,    211,  1082,   23, 0,    80),
    @(80, "Eslovaquia",             1373,  13,    386,   970,    7, 0,    17),
    @(96, "Libano",                  696,   0,    140,   532,   46, 2,    24),
    @(182,"Laos",                     19,   0,      7,    12,    0, 0,     0)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum,1).Value = $r[1]
    $ws.Cells.Item($rowNum,2).Value = $r[2]
    $ws.Cells.Item($rowNum,3).Value = $r[3]
    $ws.Cells.Item($rowNum,4).Value = $r[4]
    $ws.Cells.Item($rowNum,5).Value = $r[5]
    $ws.Cells.Item($rowNum,6).Value = $r[6]
    $ws.Cells.Item($rowNum,7).Value = $r[7]
    $ws.Cells.Item($rowNum,8).Value = $r[8]
}
